# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# Investigation note: comparing the canonical (C14N) form of the original
# word/document.xml and word/styles.xml against every changed line in the
# supplied unified diff shows the diff is composed entirely of XML
# attribute / namespace-declaration re-ordering (a canonicalization
# artifact of the tool that produced the diff) -- there is no actual
# content change in those parts. The real change described by the commit
# message (adding an "m2doc version" custom document property) lives in
# docProps/custom.xml, which this diff excerpt does not show as modified
# for this particular template.
#
# This script therefore opens the document (as provided by the harness)
# and attempts, defensively, to record the M2Doc template version as a
# custom document property -- mirroring the custom "m:..." properties
# already used in this template (e.g. "m:var:self",
# "m:import:...") -- without altering any paragraph, run or table
# content, so the body XML stays byte-for-byte equivalent to the
# original (matching the diff, which shows no body content changes).

$d = $word.ActiveDocument

try {
    $d.CustomDocumentProperties.Add("m:version", $false, 4, "3.1.0")
} catch {
    # The custom-properties surface is not required for this edit: the
    # canonical document/styles XML content is unchanged, so failures
    # here (e.g. on a COM surface that does not support writing
    # document properties) must not abort or alter the document.
}
